# Update rows 2-6 (Sending cluster: ECs, replacing prior FAPs-only rows with
# recalculated NATMI values) and append new rows 7-11 (Sending cluster: FAPs,
# recalculated values) per "Natmi following Dr Hou advice".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  ,@("ECs", "Wnt5a", "Fzd7", "ECs", 1, 0.3333333333333333, 0.1599003333333333, 0.479701, 0.0264777194346773, 0.02647771943467731, 3, 1, 1.285689666666667, 3.857069, 0.05417001519008476, 0.05417001519008476, 0.2055822062632222, 1.850239856369, 0.001434298463975272, 0.001434298463975272)
  ,@("ECs", "Wnt5a", "Fzd7", "FAPs", 1, 0.3333333333333333, 0.1599003333333333, 0.479701, 0.0264777194346773, 0.02647771943467731, 3, 1, 7.890688, 23.672064, 0.3324586794948855, 0.3324586794948855, 1.261723641429333, 11.355512772864, 0.008802747639288884, 0.008802747639288884)
  ,@("ECs", "Wnt5a", "Fzd7", "M1", 1, 0.3333333333333333, 0.1599003333333333, 0.479701, 0.0264777194346773, 0.02647771943467731, 3, 1, 6.536593666666666, 19.609781, 0.2754065676927831, 0.2754065676927832, 1.045203506164555, 9.406831555480998, 0.007292137829836974, 0.007292137829836977)
  ,@("ECs", "Wnt5a", "Fzd7", "M2", 1, 0.3333333333333333, 0.1599003333333333, 0.479701, 0.0264777194346773, 0.02647771943467731, 3, 1, 5.505106, 16.515318, 0.231946855741777, 0.2319468557417771, 0.8802682844353334, 7.922414559918, 0.006141423770086343, 0.006141423770086345)
  ,@("ECs", "Wnt5a", "Fzd7", "sCs", 1, 0.3333333333333333, 0.1599003333333333, 0.479701, 0.0264777194346773, 0.02647771943467731, 3, 1, 2.516264666666667, 7.548794, 0.1060178818804695, 0.1060178818804695, 0.4023515589548889, 3.621164030594, 0.002807111731489831, 0.002807111731489831)
  ,@("FAPs", "Wnt5a", "Fzd7", "ECs", 3, 1, 5.879152, 17.637456, 0.9735222805653226, 0.9735222805653228, 3, 1, 1.285689666666667, 3.857069, 0.05417001519008476, 0.05417001519008476, 7.558764975162667, 68.028884776464, 0.05273571672610948, 0.05273571672610949)
  ,@("FAPs", "Wnt5a", "Fzd7", "FAPs", 3, 1, 5.879152, 17.637456, 0.9735222805653226, 0.9735222805653228, 3, 1, 7.890688, 23.672064, 0.3324586794948855, 0.3324586794948855, 46.390554136576, 417.514987229184, 0.3236559318555966, 0.3236559318555967)
  ,@("FAPs", "Wnt5a", "Fzd7", "M1", 3, 1, 5.879152, 17.637456, 0.9735222805653226, 0.9735222805653228, 3, 1, 6.536593666666666, 19.609781, 0.2754065676927831, 0.2754065676927832, 38.42962772857067, 345.866649557136, 0.2681144298629461, 0.2681144298629462)
  ,@("FAPs", "Wnt5a", "Fzd7", "M2", 3, 1, 5.879152, 17.637456, 0.9735222805653226, 0.9735222805653228, 3, 1, 5.505106, 16.515318, 0.231946855741777, 0.2319468557417771, 32.36535495011201, 291.288194551008, 0.2258054319716907, 0.2258054319716907)
  ,@("FAPs", "Wnt5a", "Fzd7", "sCs", 3, 1, 5.879152, 17.637456, 0.9735222805653226, 0.9735222805653228, 3, 1, 2.516264666666667, 7.548794, 0.1060178818804695, 0.1060178818804695, 14.79350244756267, 133.141522028064, 0.1032107701489797, 0.1032107701489797)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $startRow + $i
  $rowValues = $data[$i]
  for ($c = 0; $c -lt $rowValues.Count; $c++) {
    $ws.Cells.Item($row, $c + 1).Value = $rowValues[$c]
  }
}
